# Auto-generated edit script: updates market-data cached values
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across all
# sheets, reflecting a refreshed data pull. Generated from the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1749164.5
$ws.Range("I12").Value = 4545808
$ws.Range("J12").Value = 1262.25
$ws.Range("K12").Value = 4545808
$ws.Range("L12").Value = 1262.25
$ws.Range("M12").Value = -4545638
$ws.Range("N12").Value = -1602.25
$ws.Range("H94").Value = 486.25
$ws.Range("I94").Value = 381.66666
$ws.Range("K94").Value = 381.66666
$ws.Range("M94").Value = 69.33334000000002
$ws.Range("H98").Value = 1922.4474
$ws.Range("I98").Value = 1941.4667
$ws.Range("K98").Value = 1941.4667
$ws.Range("M98").Value = -443.4666999999999
$ws.Range("H122").Value = 1922.4474
$ws.Range("I122").Value = 1941.4667
$ws.Range("K122").Value = 5824.4001
$ws.Range("M122").Value = -3374.4001
$ws.Range("H137").Value = 49787.676
$ws.Range("I137").Value = 82059.73
$ws.Range("K137").Value = 246179.19
$ws.Range("M137").Value = -243629.19
$ws.Range("H138").Value = 3981.6904
$ws.Range("J138").Value = 3915.775
$ws.Range("L138").Value = 11747.325
$ws.Range("N138").Value = -22027.325
$ws.Range("H141").Value = 12197.97
$ws.Range("I141").Value = 6690.3335
$ws.Range("K141").Value = 20071.0005
$ws.Range("M141").Value = -14891.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 4450
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1213
$ws.Range("H22").Value = 4455.091
$ws.Range("I22").Value = 1858
$ws.Range("J22").Value = 9000
$ws.Range("K22").Value = 1858
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = -1559
$ws.Range("N22").Value = -9598
$ws.Range("H32").Value = 8323.6
$ws.Range("I32").Value = 4969.5
$ws.Range("J32").Value = 24535.084
$ws.Range("K32").Value = 4969.5
$ws.Range("L32").Value = 24535.084
$ws.Range("M32").Value = -4682.5
$ws.Range("N32").Value = -25109.084
$ws.Range("H61").Value = 8181.7144
$ws.Range("I61").Value = 9655.4
$ws.Range("K61").Value = 9655.4
$ws.Range("M61").Value = -9443.4
$ws.Range("H74").Value = 23389.955
$ws.Range("I74").Value = 3458.9697
$ws.Range("K74").Value = 3458.9697
$ws.Range("M74").Value = -2584.9697
$ws.Range("H77").Value = 23389.955
$ws.Range("I77").Value = 3458.9697
$ws.Range("K77").Value = 17294.8485
$ws.Range("M77").Value = -12926.8485
$ws.Range("H92").Value = 31578
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H97").Value = 2942228
$ws.Range("J97").Value = 3849.5
$ws.Range("L97").Value = 3849.5
$ws.Range("N97").Value = -4841.5
$ws.Range("H136").Value = 8181.7144
$ws.Range("I136").Value = 9655.4
$ws.Range("K136").Value = 28966.2
$ws.Range("M136").Value = -26416.2
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 3050.25
$ws.Range("I24").Value = 1001
$ws.Range("K24").Value = 1001
$ws.Range("M24").Value = -766
$ws.Range("H25").Value = 1592.1428
$ws.Range("I25").Value = 949.4
$ws.Range("K25").Value = 949.4
$ws.Range("M25").Value = -714.4
$ws.Range("H94").Value = 3460882.8
$ws.Range("I94").Value = 4170975
$ws.Range("J94").Value = 52440
$ws.Range("K94").Value = 4170975
$ws.Range("L94").Value = 52440
$ws.Range("M94").Value = -4170524
$ws.Range("N94").Value = -53342
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19828.578
$ws.Range("I31").Value = 2457.375
$ws.Range("K31").Value = 2457.375
$ws.Range("M31").Value = -2162.375
$ws.Range("H32").Value = 8001.6665
$ws.Range("I32").Value = 6003.3335
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 6003.3335
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -5687.3335
$ws.Range("N32").Value = -10632
$ws.Range("H34").Value = 19828.578
$ws.Range("I34").Value = 2457.375
$ws.Range("K34").Value = 2457.375
$ws.Range("M34").Value = -2255.375
$ws.Range("H58").Value = 7505.0835
$ws.Range("I58").Value = 8696.4375
$ws.Range("K58").Value = 8696.4375
$ws.Range("M58").Value = -8493.4375
$ws.Range("H59").Value = 36750
$ws.Range("J59").Value = 34583.332
$ws.Range("L59").Value = 34583.332
$ws.Range("N59").Value = -36873.332
$ws.Range("H105").Value = 1477.8182
$ws.Range("I105").Value = 876.1667
$ws.Range("J105").Value = 2199.8
$ws.Range("K105").Value = 876.1667
$ws.Range("L105").Value = 2199.8
$ws.Range("M105").Value = 870.8333
$ws.Range("N105").Value = -5693.8
$ws.Range("H132").Value = 81121.664
$ws.Range("I132").Value = 60476
$ws.Range("K132").Value = 181428
$ws.Range("M132").Value = -178898
$ws.Range("H136").Value = 7505.0835
$ws.Range("I136").Value = 8696.4375
$ws.Range("K136").Value = 26089.3125
$ws.Range("M136").Value = -23539.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 467323.44
$ws.Range("I2").Value = 528.4783
$ws.Range("J2").Value = 1032391
$ws.Range("K2").Value = 3170.8698
$ws.Range("L2").Value = 6194346
$ws.Range("M2").Value = -3057.8698
$ws.Range("N2").Value = -6194572
$ws.Range("H7").Value = 73665.42999999999
$ws.Range("I7").Value = 85109.5
$ws.Range("K7").Value = 255328.5
$ws.Range("M7").Value = -255216.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 22999.5
$ws.Range("I18").Value = 10000
$ws.Range("K18").Value = 10000
$ws.Range("M18").Value = -9707
$ws.Range("H80").Value = 18730422
$ws.Range("I80").Value = 21851338
$ws.Range("J80").Value = 4924
$ws.Range("K80").Value = 21851338
$ws.Range("L80").Value = 4924
$ws.Range("M80").Value = -21850340
$ws.Range("N80").Value = -6920
$ws.Range("H83").Value = 18730422
$ws.Range("I83").Value = 21851338
$ws.Range("J83").Value = 4924
$ws.Range("K83").Value = 109256690
$ws.Range("L83").Value = 24620
$ws.Range("M83").Value = -109251698
$ws.Range("N83").Value = -34604
$ws.Range("H103").Value = 80501.60000000001
$ws.Range("J103").Value = 80501.60000000001
$ws.Range("L103").Value = 80501.60000000001
$ws.Range("N103").Value = -82845.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1743919.1
$ws.Range("I46").Value = 21739404
$ws.Range("K46").Value = 21739404
$ws.Range("M46").Value = -21739216
$ws.Range("H55").Value = 960.3333
$ws.Range("I55").Value = 1141.3334
$ws.Range("J55").Value = 688.8333
$ws.Range("K55").Value = 1141.3334
$ws.Range("L55").Value = 688.8333
$ws.Range("M55").Value = -968.3334
$ws.Range("N55").Value = -1034.8333
$ws.Range("H136").Value = 38194.05
$ws.Range("I136").Value = 53638.54
$ws.Range("K136").Value = 160915.62
$ws.Range("M136").Value = -158365.62

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 58000
$ws.Range("J82").Value = 58000
$ws.Range("L82").Value = 58000
$ws.Range("N82").Value = -58766
$ws.Range("H85").Value = 58000
$ws.Range("J85").Value = 58000
$ws.Range("L85").Value = 58000
$ws.Range("N85").Value = -60652
$ws.Range("H96").Value = 3176.3
$ws.Range("I96").Value = 2996.5
$ws.Range("J96").Value = 3446
$ws.Range("K96").Value = 2996.5
$ws.Range("L96").Value = 3446
$ws.Range("M96").Value = -1623.5
$ws.Range("N96").Value = -6192
$ws.Range("H100").Value = 2045.7
$ws.Range("J100").Value = 655.6
$ws.Range("L100").Value = 1311.2
$ws.Range("N100").Value = -2393.2
$ws.Range("H113").Value = 1567.963
$ws.Range("I113").Value = 746.9375
$ws.Range("J113").Value = 2762.182
$ws.Range("K113").Value = 2240.8125
$ws.Range("L113").Value = 8286.545999999998
$ws.Range("M113").Value = -70.8125
$ws.Range("N113").Value = -12626.546
$ws.Range("H126").Value = 1913.2632
$ws.Range("I126").Value = 2122.4285
$ws.Range("J126").Value = 1327.6
$ws.Range("K126").Value = 6367.2855
$ws.Range("L126").Value = 3982.8
$ws.Range("M126").Value = -3897.2855
$ws.Range("N126").Value = -8922.799999999999
$ws.Range("H136").Value = 2819.575
$ws.Range("I136").Value = 2679.257
$ws.Range("J136").Value = 3801.8
$ws.Range("K136").Value = 8037.771000000001
$ws.Range("L136").Value = 11405.4
$ws.Range("M136").Value = -5487.771000000001
$ws.Range("N136").Value = -16505.4
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
